# Update Betfair Back/Lay odds values for 2025-11-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Deportes Concepcion x Antofagasta)
$ws.Range("F2").Value = 2.08
$ws.Range("G2").Value = 2.48
$ws.Range("H2").Value = 3.45
$ws.Range("I2").Value = 4.7
$ws.Range("K2").Value = 4.3
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.73
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.27
$ws.Range("S2").Value = 3.45
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 1.94
$ws.Range("V2").Value = 1.27
$ws.Range("W2").Value = 1.67

# Row 3 (Botafogo FR x Sport Recife)
$ws.Range("F3").Value = 1.27
$ws.Range("I3").Value = 15.5
$ws.Range("N3").Value = 5.1
$ws.Range("Q3").Value = 1.63
$ws.Range("AF3").Value = 8.199999999999999

# Row 4
$ws.Range("AK4").Value = 29

# Row 6
$ws.Range("I6").Value = 26
$ws.Range("J6").Value = 6.8

# Row 7
$ws.Range("F7").Value = 1.28
$ws.Range("H7").Value = 10.5
$ws.Range("J7").Value = 5.4
$ws.Range("O7").Value = 1.26
$ws.Range("P7").Value = 2.1
$ws.Range("AB7").Value = 9
$ws.Range("AD7").Value = 55
$ws.Range("AK7").Value = 18

# Row 8
$ws.Range("F8").Value = 1.86
$ws.Range("G8").Value = 1.98
$ws.Range("H8").Value = 4.7
$ws.Range("I8").Value = 5.3
$ws.Range("K8").Value = 3.9
$ws.Range("T8").Value = 1.96
$ws.Range("U8").Value = 1.85
$ws.Range("V8").Value = 1.23
$ws.Range("W8").Value = 2.02
$ws.Range("AD8").Value = 21
$ws.Range("AI8").Value = 85
$ws.Range("AJ8").Value = 22
$ws.Range("AN8").Value = 17

# Row 9
$ws.Range("X9").Value = 11

# Row 10
$ws.Range("G10").Value = 1.11
$ws.Range("H10").Value = 30
$ws.Range("J10").Value = 13.5
$ws.Range("N10").Value = 9
$ws.Range("O10").Value = 1.09
$ws.Range("Q10").Value = 1.29
$ws.Range("R10").Value = 2.08
$ws.Range("S10").Value = 1.74
$ws.Range("T10").Value = 2.54
$ws.Range("U10").Value = 1.52
$ws.Range("W10").Value = 9.800000000000001
$ws.Range("X10").Value = 60
$ws.Range("Y10").Value = 170
$ws.Range("AB10").Value = 15.5
$ws.Range("AC10").Value = 980
$ws.Range("AD10").Value = 160
$ws.Range("AF10").Value = 9.6
$ws.Range("AG10").Value = 18.5
$ws.Range("AH10").Value = 80
$ws.Range("AJ10").Value = 8.199999999999999
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 70
$ws.Range("AM10").Value = 510
$ws.Range("AN10").Value = 2.4
